$wb = $excel.ActiveWorkbook

# --- "Channel Map" sheet: split the old "600 700" / "abs" row into two rows,
#     one per wavelength (600 -> abs600, 700 -> abs700) ---
$wsChan = $wb.Worksheets.Item("Channel Map")
$wsChan.Range("A2").Value = "600"
$wsChan.Range("B2").Value = "abs600"
$wsChan.Range("A3").Value = "700"
$wsChan.Range("B3").Value = "abs700"
# keep the Channel column as text so "600"/"700" aren't read back as numbers
$wsChan.Range("A2:A3").NumberFormat = "@"

# --- "Samples" sheet: fix the Channels value on the existing row and add a
#     second spectramax sample row ---
$wsSamples = $wb.Worksheets.Item("Samples")
$wsSamples.Range("C2").Value = "535_485,600,700"

$wsSamples.Range("A3").Value = "Plate reader"
$wsSamples.Range("B3").Value = "`$GITHUB_WORKSPACE/test/inputs/spectramax-data2.txt"
$wsSamples.Range("C3").Value = "600,700,530_485_1,530_485_2,530_485_3"
$wsSamples.Range("D3").Value = "spectramax"
$wsSamples.Range("E3").Value = 2

# Channels column holds comma separated text - format as text
$wsSamples.Range("C1:C3").NumberFormat = "@"

# match the new column C width to the other bestFit columns
$wsSamples.Columns.Item(3).ColumnWidth = 10.330729166666666

# --- refresh the selections shown when the sheets were last viewed ---
# select on the non-active sheet first, then the active sheet last so the
# active sheet stays the selected tab
$wsChan.Range("B4").Select()
$wsSamples.Range("C4").Select()
